$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 105
$ws1.Range("F3").Value  = 12086
$ws1.Range("F4").Value  = 46
$ws1.Range("F8").Value  = 11982
$ws1.Range("F9").Value  = 503
$ws1.Range("F13").Value = 1800
$ws1.Range("F14").Value = 5926
$ws1.Range("F15").Value = 132
$ws1.Range("F16").Value = 3558

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 579
$ws2.Range("F4").Value = 10

# Sheet 4: 全部类型 (All Types - combined view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 579
$ws4.Range("F3").Value  = 105
$ws4.Range("F5").Value  = 12086
$ws4.Range("F6").Value  = 46
$ws4.Range("F8").Value  = 10
$ws4.Range("F11").Value = 11982
$ws4.Range("F12").Value = 503
$ws4.Range("F16").Value = 1800
$ws4.Range("F18").Value = 5926
$ws4.Range("F19").Value = 132
$ws4.Range("F20").Value = 3558
